$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Metadata")
$ws2 = $wb.Worksheets.Item("Include from Engagement Commu")

# Version: 5.0.0 -> 6.0.0 (row 3)
$ws1.Range("B3").Value = "6.0.0"

# Date (row 8)
$ws1.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher row (row 9): B9 empty -> "Alvearie Team"
$ws1.Range("B9").Value = "Alvearie Team"

# Old row 10/11 were "Contact" / "No display for ContactDetail" (duplicated).
# New row 10 becomes Jurisdiction / United States of America
$ws1.Range("A10").Value = "Jurisdiction"
$ws1.Range("B10").Value = "United States of America"

# New row 11 becomes Description / Communication modes supported for Engagement communication
$ws1.Range("A11").Value = "Description"
$ws1.Range("B11").Value = "Communication modes supported for Engagement communication"

# New row 12 becomes Purpose / empty
$ws1.Range("A12").Value = "Purpose"
$ws1.Range("B12").ClearContents()

# New row 13 becomes Copyright / empty
$ws1.Range("A13").Value = "Copyright"
$ws1.Range("B13").ClearContents()

# New row 14 becomes Immutable / BooleanType[null]
$ws1.Range("A14").Value = "Immutable"
$ws1.Range("B14").Value = "BooleanType[null]"

# Delete old row 15 (Immutable/BooleanType[null]) which is now redundant since rows shifted up.
$ws1.Rows.Item(15).Delete()
